# Update the exposure-site record from "Pascoe Vale / Elite Swimming" to
# "Point Cook / The Coffeeologist Cafe", fix up the exposure-period times,
# and flip the Exist flags (row 2 -> old, row 3 -> new).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A3").Value = "Point Cook"
$ws.Range("B2:B3").Value = "The Coffeeologist Cafe  70/300 Point Cook Rd  Point Cook VIC 3030"
$ws.Range("C2").Value = "11:30am - 12:10pm  0/2/2021"
$ws.Range("C3").Value = "11:30am - 12:10pm  10/2/2021"
$ws.Range("E2").Value = "old"
$ws.Range("E3").Value = "new"

# Select the edited columns and let them re-fit their width to the new,
# longer text (mirrors the wider A:C columns seen after the edit).
[void]$ws.Columns("A:E").Select()
[void]$ws.Columns("A:C").AutoFit()

$wb.Save() | Out-Null
